$wb = $excel.ActiveWorkbook

# --- Rushing sheet ---
$ws = $wb.Worksheets.Item("Rushing")

# Row 2: P.Mahomes
$ws.Range("C2").Value = 18
$ws.Range("D2").Value = 20
$ws.Range("E2").Value = 12
$ws.Range("F2").Value = 15

# Row 4: C.Edwards-Helaire
$ws.Range("C4").Value = 77
$ws.Range("D4").Value = 47
$ws.Range("F4").Value = 15

# Row 6: J.McKinnon
$ws.Range("C6").Value = 19
$ws.Range("D6").Value = 9
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 6

# Row 9: T.Hill
$ws.Range("C9").Value = 8

# Row 10: M.Hardman
$ws.Range("C10").Value = 10

# Row 12: B.Bell
$ws.Range("E12").Value = 5
$ws.Range("F12").Value = 2

# --- Receiving sheet ---
$ws2 = $wb.Worksheets.Item("Receiving")

# Row 2: C.Edwards-Helaire
$ws2.Range("C2").Value = 39
$ws2.Range("D2").Value = 30
$ws2.Range("G2").Value = 5

# Row 4: J.McKinnon
$ws2.Range("C4").Value = 15
$ws2.Range("D4").Value = 12
$ws2.Range("E4").Value = 4
$ws2.Range("F4").Value = 3
$ws2.Range("G4").Value = 5
$ws2.Range("H4").Value = 4

# Row 7: T.Hill
$ws2.Range("C7").Value = 134
$ws2.Range("D7").Value = 107
$ws2.Range("E7").Value = 40
$ws2.Range("F7").Value = 17
$ws2.Range("G7").Value = 28

# Row 8: M.Hardman
$ws2.Range("E8").Value = 19
$ws2.Range("F8").Value = 10

# Row 9: B.Pringle
$ws2.Range("C9").Value = 51
$ws2.Range("D9").Value = 35
$ws2.Range("G9").Value = 7
$ws2.Range("H9").Value = 5

# Row 10: D.Robinson
$ws2.Range("C10").Value = 32

# Row 13: T.Kelce
$ws2.Range("C13").Value = 125
$ws2.Range("D13").Value = 89
$ws2.Range("E13").Value = 27
$ws2.Range("F13").Value = 17
$ws2.Range("G13").Value = 21
$ws2.Range("H13").Value = 17

# Row 15: N.Gray
$ws2.Range("C15").Value = 11
$ws2.Range("D15").Value = 6
